# Spell-check fix: the "_GoBack" last-edit bookmark that Word stamps
# at the point of the most recent edit needs to sit inside the word
# "originalite" (between "orig" and "inalite") instead of after the
# closing parenthesis at the end of the paragraph.
#
# We locate "inalite" (unique in the document) via Find, collapse a
# Range to its start, and re-Add the "_GoBack" bookmark there. Word
# bookmark names are unique, so Bookmarks.Add with the existing name
# "_GoBack" simply relocates it (and splits the run) rather than
# creating a duplicate.

$d = $word.ActiveDocument

$findRange = $d.Content
$findRange.Find.Execute("inalite", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

$pos = $findRange.Start
$target = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $target)
